$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing header (beta_152 -> PL_beta) and add new headers ---
$ws.Range("AT1").Value = "PL_beta"
$ws.Range("AZ1").Value = "WTAR_Standard"
$ws.Range('BA1').Value = "ACE_(Attention/18)'"
$ws.Range('BB1').Value = "ACE_(Memory/26)'"
$ws.Range("AZ1:BB1").NumberFormat = "0.00"

# --- Fill in new data columns: AZ=WTAR_Standard, BA=ACE_(Attention/18), BB=ACE_(Memory/26) ---
$ws.Range("AZ2").Value = 119
$ws.Range("AZ2").NumberFormat = "0.00"
$ws.Range("AZ3").Value = 127
$ws.Range("AZ3").NumberFormat = "0.00"
$ws.Range("AZ4").Value = 117
$ws.Range("AZ4").NumberFormat = "0.00"
$ws.Range("AZ5").Value = 115
$ws.Range("AZ5").NumberFormat = "0.00"
$ws.Range("AZ6").Value = 113
$ws.Range("AZ6").NumberFormat = "0.00"
$ws.Range("AZ7").Value = 103
$ws.Range("AZ7").NumberFormat = "0.00"
$ws.Range("AZ8").Value = 115
$ws.Range("AZ8").NumberFormat = "0.00"
$ws.Range("AZ9").Value = 108
$ws.Range("AZ9").NumberFormat = "0.00"
$ws.Range("AZ10").Value = 129
$ws.Range("AZ10").NumberFormat = "0.00"
$ws.Range("AZ11").Value = 119
$ws.Range("AZ11").NumberFormat = "0.00"
$ws.Range("AZ12").Value = 111
$ws.Range("AZ12").NumberFormat = "0.00"
$ws.Range("AZ13").Value = 117
$ws.Range("AZ13").NumberFormat = "0.00"
$ws.Range("AZ14").Value = 129
$ws.Range("AZ14").NumberFormat = "0.00"
$ws.Range("AZ15").Value = 99
$ws.Range("AZ15").NumberFormat = "0.00"
$ws.Range("AZ16").Value = 101
$ws.Range("AZ16").NumberFormat = "0.00"
$ws.Range("AZ17").Value = 129
$ws.Range("AZ17").NumberFormat = "0.00"
$ws.Range("AZ18").Value = 106
$ws.Range("AZ18").NumberFormat = "0.00"
$ws.Range("AZ19").Value = 117
$ws.Range("AZ19").NumberFormat = "0.00"
$ws.Range("AZ20").Value = 122
$ws.Range("AZ20").NumberFormat = "0.00"
$ws.Range("AZ21").Value = 104
$ws.Range("AZ21").NumberFormat = "0.00"
$ws.Range("AZ22").Value = 119
$ws.Range("AZ22").NumberFormat = "0.00"
$ws.Range("AZ23").Value = 122
$ws.Range("AZ23").NumberFormat = "0.00"
$ws.Range("AZ24").Value = 113
$ws.Range("AZ24").NumberFormat = "0.00"
$ws.Range("AZ25").Value = 113
$ws.Range("AZ25").NumberFormat = "0.00"
$ws.Range("AZ26").Value = 115
$ws.Range("AZ26").NumberFormat = "0.00"
$ws.Range("AZ27").Value = 117
$ws.Range("AZ27").NumberFormat = "0.00"
$ws.Range("AZ28").Value = 122
$ws.Range("AZ28").NumberFormat = "0.00"
$ws.Range("AZ29").Value = 111
$ws.Range("AZ29").NumberFormat = "0.00"
$ws.Range("AZ30").Value = 127
$ws.Range("AZ30").NumberFormat = "0.00"
$ws.Range("AZ31").Value = 124
$ws.Range("AZ31").NumberFormat = "0.00"
$ws.Range("AZ32").Value = 126
$ws.Range("AZ32").NumberFormat = "0.00"
$ws.Range("AZ33").Value = 122
$ws.Range("AZ33").NumberFormat = "0.00"
$ws.Range("AZ34").Value = 124
$ws.Range("AZ34").NumberFormat = "0.00"
$ws.Range("AZ35").Value = 101
$ws.Range("AZ35").NumberFormat = "0.00"
$ws.Range("AZ36").Value = 126
$ws.Range("AZ36").NumberFormat = "0.00"
$ws.Range("AZ37").Value = 119
$ws.Range("AZ37").NumberFormat = "0.00"
$ws.Range("AZ38").Value = 122
$ws.Range("AZ38").NumberFormat = "0.00"
$ws.Range("AZ39").Value = 122
$ws.Range("AZ39").NumberFormat = "0.00"
$ws.Range("AZ40").Value = 124
$ws.Range("AZ40").NumberFormat = "0.00"
$ws.Range("AZ41").Value = 122
$ws.Range("AZ41").NumberFormat = "0.00"
$ws.Range("AZ42").Value = 122
$ws.Range("AZ42").NumberFormat = "0.00"
$ws.Range("AZ43").Value = 127
$ws.Range("AZ43").NumberFormat = "0.00"
$ws.Range("AZ44").Value = 106
$ws.Range("AZ44").NumberFormat = "0.00"
$ws.Range("AZ45").Value = 126
$ws.Range("AZ45").NumberFormat = "0.00"
$ws.Range("AZ46").Value = 120
$ws.Range("AZ46").NumberFormat = "0.00"
$ws.Range("AZ47").Value = 127
$ws.Range("AZ47").NumberFormat = "0.00"
$ws.Range("AZ48").Value = 122
$ws.Range("AZ48").NumberFormat = "0.00"
$ws.Range("AZ49").Value = 122
$ws.Range("AZ49").NumberFormat = "0.00"
$ws.Range("AZ50").Value = 101
$ws.Range("AZ50").NumberFormat = "0.00"
$ws.Range("AZ51").Value = 117
$ws.Range("AZ51").NumberFormat = "0.00"
$ws.Range("AZ52").Value = 122
$ws.Range("AZ52").NumberFormat = "0.00"
$ws.Range("AZ53").Value = 119
$ws.Range("AZ53").NumberFormat = "0.00"
$ws.Range("AZ54").Value = 120
$ws.Range("AZ54").NumberFormat = "0.00"
$ws.Range("AZ55").Value = 120
$ws.Range("AZ55").NumberFormat = "0.00"
$ws.Range("AZ56").Value = 126
$ws.Range("AZ56").NumberFormat = "0.00"
$ws.Range("AZ57").Value = 111
$ws.Range("AZ57").NumberFormat = "0.00"
$ws.Range("AZ58").Value = 127
$ws.Range("AZ58").NumberFormat = "0.00"
$ws.Range("AZ59").Value = 117
$ws.Range("AZ59").NumberFormat = "0.00"
$ws.Range("AZ60").Value = 126
$ws.Range("AZ60").NumberFormat = "0.00"
$ws.Range("AZ61").Value = 120
$ws.Range("AZ61").NumberFormat = "0.00"
$ws.Range("AZ62").Value = 124
$ws.Range("AZ62").NumberFormat = "0.00"
$ws.Range("AZ63").Value = 122
$ws.Range("AZ63").NumberFormat = "0.00"
$ws.Range("AZ64").Value = 106
$ws.Range("AZ64").NumberFormat = "0.00"
$ws.Range("AZ65").Value = 119
$ws.Range("AZ65").NumberFormat = "0.00"
$ws.Range("AZ66").Value = 111
$ws.Range("AZ66").NumberFormat = "0.00"
$ws.Range("AZ67").Value = 119
$ws.Range("AZ67").NumberFormat = "0.00"
$ws.Range("AZ68").Value = 108
$ws.Range("AZ68").NumberFormat = "0.00"
$ws.Range("AZ69").Value = 113
$ws.Range("AZ69").NumberFormat = "0.00"
$ws.Range("AZ70").Value = 122
$ws.Range("AZ70").NumberFormat = "0.00"
$ws.Range("AZ71").Value = 122
$ws.Range("AZ71").NumberFormat = "0.00"
$ws.Range("AZ72").Value = 122
$ws.Range("AZ72").NumberFormat = "0.00"
$ws.Range("AZ73").Value = 124
$ws.Range("AZ73").NumberFormat = "0.00"
$ws.Range("AZ74").Value = 124
$ws.Range("AZ74").NumberFormat = "0.00"
$ws.Range("AZ75").Value = 122
$ws.Range("AZ75").NumberFormat = "0.00"
$ws.Range("AZ76").Value = 124
$ws.Range("AZ76").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 122; $arr[0,1] = 14; $arr[0,2] = 24
$ws.Range("AZ77:BB77").Value = $arr
$ws.Range("AZ77:BB77").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 24
$ws.Range("AZ78:BB78").Value = $arr
$ws.Range("AZ78:BB78").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ79:BB79").Value = $arr
$ws.Range("AZ79:BB79").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ80:BB80").Value = $arr
$ws.Range("AZ80:BB80").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 17; $arr[0,2] = 23
$ws.Range("AZ81:BB81").Value = $arr
$ws.Range("AZ81:BB81").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 14; $arr[0,2] = 24
$ws.Range("AZ82:BB82").Value = $arr
$ws.Range("AZ82:BB82").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ83:BB83").Value = $arr
$ws.Range("AZ83:BB83").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 93; $arr[0,1] = 17; $arr[0,2] = 22
$ws.Range("AZ84:BB84").Value = $arr
$ws.Range("AZ84:BB84").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 110; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ85:BB85").Value = $arr
$ws.Range("AZ85:BB85").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ86:BB86").Value = $arr
$ws.Range("AZ86:BB86").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ87:BB87").Value = $arr
$ws.Range("AZ87:BB87").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ88:BB88").Value = $arr
$ws.Range("AZ88:BB88").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 113; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ89:BB89").Value = $arr
$ws.Range("AZ89:BB89").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 16; $arr[0,2] = 26
$ws.Range("AZ90:BB90").Value = $arr
$ws.Range("AZ90:BB90").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 111; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ91:BB91").Value = $arr
$ws.Range("AZ91:BB91").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ92:BB92").Value = $arr
$ws.Range("AZ92:BB92").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ93:BB93").Value = $arr
$ws.Range("AZ93:BB93").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ94:BB94").Value = $arr
$ws.Range("AZ94:BB94").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 113; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ95:BB95").Value = $arr
$ws.Range("AZ95:BB95").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ96:BB96").Value = $arr
$ws.Range("AZ96:BB96").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ97:BB97").Value = $arr
$ws.Range("AZ97:BB97").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 118; $arr[0,1] = 17; $arr[0,2] = 17
$ws.Range("AZ98:BB98").Value = $arr
$ws.Range("AZ98:BB98").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 15; $arr[0,2] = 26
$ws.Range("AZ99:BB99").Value = $arr
$ws.Range("AZ99:BB99").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 114; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ100:BB100").Value = $arr
$ws.Range("AZ100:BB100").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 16; $arr[0,2] = 26
$ws.Range("AZ101:BB101").Value = $arr
$ws.Range("AZ101:BB101").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 118; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ102:BB102").Value = $arr
$ws.Range("AZ102:BB102").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ103:BB103").Value = $arr
$ws.Range("AZ103:BB103").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 114; $arr[0,1] = 18; $arr[0,2] = 21
$ws.Range("AZ104:BB104").Value = $arr
$ws.Range("AZ104:BB104").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ105:BB105").Value = $arr
$ws.Range("AZ105:BB105").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ106:BB106").Value = $arr
$ws.Range("AZ106:BB106").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ107:BB107").Value = $arr
$ws.Range("AZ107:BB107").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 112; $arr[0,1] = 16; $arr[0,2] = 22
$ws.Range("AZ108:BB108").Value = $arr
$ws.Range("AZ108:BB108").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 116; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ109:BB109").Value = $arr
$ws.Range("AZ109:BB109").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ110:BB110").Value = $arr
$ws.Range("AZ110:BB110").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ111:BB111").Value = $arr
$ws.Range("AZ111:BB111").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 101; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ112:BB112").Value = $arr
$ws.Range("AZ112:BB112").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 118; $arr[0,1] = 16; $arr[0,2] = 26
$ws.Range("AZ113:BB113").Value = $arr
$ws.Range("AZ113:BB113").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 108; $arr[0,1] = 16; $arr[0,2] = 26
$ws.Range("AZ114:BB114").Value = $arr
$ws.Range("AZ114:BB114").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 17; $arr[0,2] = 26
$ws.Range("AZ115:BB115").Value = $arr
$ws.Range("AZ115:BB115").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ116:BB116").Value = $arr
$ws.Range("AZ116:BB116").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ117:BB117").Value = $arr
$ws.Range("AZ117:BB117").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 115; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ118:BB118").Value = $arr
$ws.Range("AZ118:BB118").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 119; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ119:BB119").Value = $arr
$ws.Range("AZ119:BB119").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ120:BB120").Value = $arr
$ws.Range("AZ120:BB120").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 117; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ121:BB121").Value = $arr
$ws.Range("AZ121:BB121").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ122:BB122").Value = $arr
$ws.Range("AZ122:BB122").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 119; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ123:BB123").Value = $arr
$ws.Range("AZ123:BB123").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 112; $arr[0,1] = 16; $arr[0,2] = 25
$ws.Range("AZ124:BB124").Value = $arr
$ws.Range("AZ124:BB124").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 116; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ125:BB125").Value = $arr
$ws.Range("AZ125:BB125").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ126:BB126").Value = $arr
$ws.Range("AZ126:BB126").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 116; $arr[0,1] = 18; $arr[0,2] = 24
$ws.Range("AZ127:BB127").Value = $arr
$ws.Range("AZ127:BB127").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 24
$ws.Range("AZ128:BB128").Value = $arr
$ws.Range("AZ128:BB128").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 118; $arr[0,1] = 17; $arr[0,2] = 22
$ws.Range("AZ129:BB129").Value = $arr
$ws.Range("AZ129:BB129").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 113; $arr[0,1] = 17; $arr[0,2] = 24
$ws.Range("AZ130:BB130").Value = $arr
$ws.Range("AZ130:BB130").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 17; $arr[0,2] = 23
$ws.Range("AZ131:BB131").Value = $arr
$ws.Range("AZ131:BB131").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ132:BB132").Value = $arr
$ws.Range("AZ132:BB132").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 18; $arr[0,2] = 24
$ws.Range("AZ133:BB133").Value = $arr
$ws.Range("AZ133:BB133").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 120; $arr[0,1] = 16; $arr[0,2] = 24
$ws.Range("AZ134:BB134").Value = $arr
$ws.Range("AZ134:BB134").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 80; $arr[0,1] = 11; $arr[0,2] = 20
$ws.Range("AZ135:BB135").Value = $arr
$ws.Range("AZ135:BB135").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 16; $arr[0,2] = 23
$ws.Range("AZ136:BB136").Value = $arr
$ws.Range("AZ136:BB136").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 15; $arr[0,2] = 26
$ws.Range("AZ137:BB137").Value = $arr
$ws.Range("AZ137:BB137").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 16; $arr[0,2] = 26
$ws.Range("AZ138:BB138").Value = $arr
$ws.Range("AZ138:BB138").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 15; $arr[0,2] = 26
$ws.Range("AZ139:BB139").Value = $arr
$ws.Range("AZ139:BB139").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 122; $arr[0,1] = 14; $arr[0,2] = 26
$ws.Range("AZ140:BB140").Value = $arr
$ws.Range("AZ140:BB140").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 117; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ141:BB141").Value = $arr
$ws.Range("AZ141:BB141").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ142:BB142").Value = $arr
$ws.Range("AZ142:BB142").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 116; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ143:BB143").Value = $arr
$ws.Range("AZ143:BB143").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 119; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ144:BB144").Value = $arr
$ws.Range("AZ144:BB144").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 107; $arr[0,1] = 18; $arr[0,2] = 24
$ws.Range("AZ145:BB145").Value = $arr
$ws.Range("AZ145:BB145").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 18; $arr[0,2] = 26
$ws.Range("AZ146:BB146").Value = $arr
$ws.Range("AZ146:BB146").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 128; $arr[0,1] = 18; $arr[0,2] = 19
$ws.Range("AZ147:BB147").Value = $arr
$ws.Range("AZ147:BB147").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 118; $arr[0,1] = 17; $arr[0,2] = 24
$ws.Range("AZ148:BB148").Value = $arr
$ws.Range("AZ148:BB148").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 122; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ149:BB149").Value = $arr
$ws.Range("AZ149:BB149").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 123; $arr[0,1] = 17; $arr[0,2] = 25
$ws.Range("AZ150:BB150").Value = $arr
$ws.Range("AZ150:BB150").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 112; $arr[0,1] = 18; $arr[0,2] = 23
$ws.Range("AZ151:BB151").Value = $arr
$ws.Range("AZ151:BB151").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 112; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ152:BB152").Value = $arr
$ws.Range("AZ152:BB152").NumberFormat = "0.00"
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 125; $arr[0,1] = 18; $arr[0,2] = 25
$ws.Range("AZ153:BB153").Value = $arr
$ws.Range("AZ153:BB153").NumberFormat = "0.00"

# --- Restore final selection/view state ---
$ws.Range("AT36").Select()
